$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet renders a sequence diagram. Rows 5-21 use columns B/C/D as
# "lifeline gutter" columns that host call/return arrows (and thick-left-
# border style markers) at increasing nesting depth. This edit shifts the
# nesting of the B/C/D gutters so the call/return markers line up with the
# correct depth: old column C becomes the new column B, and old column D
# becomes the new column C. Rows 8, 9 and 10 gain a brand new column D
# marker (previously they had no column D content at all).

# ---------------------------------------------------------------------
# Part 1: rows 5-11 originally only used columns A-C (no column D).
# Shift: new B = old C ; new C = removed entirely.
# ---------------------------------------------------------------------
$rowsNoD = 5,6,7,8,9,10,11
foreach ($r in $rowsNoD) {
    $ws.Range("B$r").ClearContents()
    $ws.Range("C$r").Copy($ws.Range("B$r"))
    $ws.Range("C$r").Clear()
}

# ---------------------------------------------------------------------
# Part 2: rows 12-21 originally used columns A, C and D.
# Shift: new B = old C ; new C = old D ; new D = removed entirely.
# ---------------------------------------------------------------------
$rowsWithD = 12,13,14,15,16,17,18,19,20,21
foreach ($r in $rowsWithD) {
    $ws.Range("B$r").ClearContents()
    $ws.Range("C$r").Copy($ws.Range("B$r"))
    $ws.Range("C$r").ClearContents()
    $ws.Range("D$r").Copy($ws.Range("C$r"))
    $ws.Range("D$r").Clear()
}

# ---------------------------------------------------------------------
# Part 3: rows 8 and 10 gain a brand new column-D arrow marker (call /
# return arrows that used to be drawn one gutter to the left). Row 9
# (between them) gains a matching but empty/styled column-D gutter cell
# so the vertical lifeline border is continuous.
# Style is cloned from the already-present "gutter" cell in column A of
# the same row (same thick-left-border style used throughout the gutter
# columns), then the arrow glyph is written into it.
# ---------------------------------------------------------------------
$ws.Range("A8").Copy($ws.Range("D8"))
$ws.Range("D8").Formula = "→"

$ws.Range("A9").Copy($ws.Range("D9"))

$ws.Range("A10").Copy($ws.Range("D10"))
$ws.Range("D10").Formula = "←"
